$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 274 (shifts old 274:293 down to 275:294)
$ws.Rows.Item(274).Insert()

# Populate the newly inserted row 274 with the new record
$ws.Range("A274").Value = 11
$ws.Range("B274").Value = "Vega Monumental Concepción"
$ws.Range("C274").Value = "Bíobío"
$ws.Range("D274").Value = 45265
$ws.Range("E274").Value = 8
$ws.Range("F274").Value = 100112032
$ws.Range("G274").Value = "Zapallo italiano"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 100
$ws.Range("K274").Value = 10000
$ws.Range("L274").Value = 11000
$ws.Range("M274").Value = 10500
$ws.Range("N274").Value = "$/caja 50 unidades"
$ws.Range("O274").Value = "Región de O'Higgins"
$ws.Range("P274").Value = 210
$ws.Range("Q274").Value = 50
$ws.Range("R274").Value = "Hortaliza"
